{"js": "// Word JS API (Office.js) edit script.\n// Body of: async (context) => { ... }\n//\n// 1) Turn \"GitHub repo\" into \"GitHub repository\" by placing the\n//    insertion point right after \"repo\" and typing \"sitory\" there\n//    (mirrors how the author actually fixed the readme text).\n// 2) Make the section's page orientation explicit (Portrait), which is\n//    what produced the added w:orient=\"portrait\" attribute on <w:pgSz>.\n\n// --- 1. Fix \"repo\" -> \"repository\" -------------------------------------\nconst repoResults = context.document.body.search(\"repo\", {\n  matchCase: false,\n  matchWholeWord: false\n});\nrepoResults.load(\"text\");\nawait context.sync();\n\nif (repoResults.items.length > 0) {\n  const repoRange = repoResults.items[0];\n  // Collapsed range sitting right after the \"o\" in \"repo\".\n  const insertionPoint = repoRange.getRange(Word.RangeLocation.end);\n  insertionPoint.insertText(\"sitory\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2. Make page orientation explicit (Portrait) -----------------------\nconst pageSetup = context.document.pageSetup;\npageSetup.orientation = Word.PageOrientation.portrait;\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word / $app -> Word.Application, $d -> $word.ActiveDocument\n#\n# 1) Turn \"GitHub repo\" into \"GitHub repository\" by finding \"repo\",\n#    collapsing the found range to its end, and typing \"sitory\" there\n#    (mirrors how the author actually fixed the readme text).\n# 2) Make the section's page orientation explicit (Portrait), which is\n#    what produced the added w:orient=\"portrait\" attribute on <w:pgSz>.\n\n$d = $word.ActiveDocument\n\n# --- 1. Fix \"repo\" -> \"repository\" --------------------------------------\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"repo\")\nif ($found) {\n    $findRange.Collapse(0)   # wdCollapseEnd\n    $findRange.InsertAfter(\"sitory\")\n}\n\n# --- 2. Make page orientation explicit (Portrait) ------------------------\n$d.PageSetup.Orientation = 0   # wdOrientPortrait\n"}
